$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns C, D, E, F, G, H, J, N
# (columns B, I, K, L, M, O are unchanged - left as-is)
$newValues = @{
    2 = @(0.04837263193822139, 0.08076544123092333, 0.1133092308110477, 3.954288008995889, 3.5176108167297, 2.358650277960862, 0.2323591097579509, 2.873581687696401)
    3 = @(0.04291586150900173, 0.07849806491704214, 0.1129899566129922, 3.771370646392, 3.301587764275439, 2.267373502473845, 0.2254799030834818, 2.562605684679681)
    4 = @(0.03959561195001982, 0.07708835093855981, 0.1128551022632927, 3.662213813342333, 3.17181790680678, 2.213175233925938, 0.2214329583294443, 2.371325805375818)
    5 = @(0.03824977303857224, 0.07650930721145954, 0.1128154488327873, 3.61850868487744, 3.119638749499074, 2.191544775453224, 0.2198275908243801, 2.293303068605894)
    6 = @(0.03802671988772488, 0.0764128764471792, 0.112809786151427, 3.611297953766808, 3.111016377872318, 2.187980339722969, 0.2195636460189974, 2.2803432614038)
    7 = @(0.03957743292303917, 0.07708056048665668, 0.1128545056388823, 3.661621265120232, 3.171111377513114, 2.212881683554258, 0.2214111313566605, 2.370273851392596)
    8 = @(0.04648461136100934, 0.07998718910423719, 0.1131863952433285, 3.890553302341232, 3.44251932213632, 2.326789733958549, 0.2299500163077965, 2.766433886209882)
    9 = @(0.06028899607788674, 0.0855556671229607, 0.1143267567236705, 4.365297793557602, 3.99839030280333, 2.565213172815561, 0.2481315412707659, 3.540180268007646)
    10 = @(0.07061782457429899, 0.08957788668929112, 0.1154691160534504, 4.730963310585821, 4.422482990430353, 2.750149006206982, 0.2624139554790048, 4.10623028343673)
    11 = @(0.07536367885205664, 0.09139523579088404, 0.1160563421568135, 4.90123507854571, 4.619111932801957, 2.836535916717708, 0.2691231557129896, 4.363110593465422)
    12 = @(0.07716812685862351, 0.09208184401926189, 0.116288541591306, 4.966298633669908, 4.694127045197376, 2.869584196910012, 0.2716951551983016, 4.460285735713398)
    13 = @(0.07677917397765555, 0.09193403782308707, 0.116238094238188, 4.95225966812967, 4.677946097155541, 2.862451566768527, 0.2711398195665708, 4.439361943450422)
    14 = @(0.07551198269777615, 0.09145175402413486, 0.1160752476701958, 4.906576046229361, 4.625272194154491, 2.839248033856052, 0.2693341218030696, 4.371107314139294)
    15 = @(0.07473675688960668, 0.09115614115602, 0.1159767830357055, 4.878670353067548, 4.593081008551223, 2.825079204977101, 0.268232191764298, 4.329286057409945)
    16 = @(0.07030866489481014, 0.08945888632213439, 0.1154321086230858, 4.71991643304068, 4.409709365074946, 2.744549769859248, 0.2619798349272742, 4.089429168003562)
    17 = @(0.06760462282041146, 0.08841464377878339, 0.1151153524600943, 4.623546572525612, 4.298182346519695, 2.695733525950004, 0.2581991445066052, 3.94211849063862)
    18 = @(0.06605374256338337, 0.0878128457814995, 0.1149395147330772, 4.568485704440718, 4.234382464436237, 2.667867652809321, 0.2560444970519029, 3.857331695637129)
    19 = @(0.065529383964261, 0.08760887918698756, 0.11488106665945, 4.549905725111984, 4.212839919092346, 2.658468829434128, 0.2553183630231501, 3.828614786364199)
    20 = @(0.06789201234747111, 0.08852592561328265, 0.1151484134686918, 4.633767022025438, 4.310018458435081, 2.700908087636094, 0.2585995390348472, 3.957806003277312)
    21 = @(0.07588398537247087, 0.09159345381029738, 0.1161228119448872, 4.919978381715907, 4.640728520165908, 2.846054283815135, 0.269863640468273, 4.391158149571254)
    22 = @(0.08114997452075556, 0.09358913304082961, 0.1168169890189503, 5.110455770377826, 4.860118637852167, 2.942876066604526, 0.2774086648391574, 4.673791817957863)
    23 = @(0.07833533198132159, 0.09252476986636538, 0.1164412069167753, 5.008474398812893, 4.742720673453675, 2.891017476650177, 0.27336467515957, 4.523002190001307)
    24 = @(0.06776207202669582, 0.08847561962414829, 0.1151334470668743, 4.629145288575558, 4.304666361949273, 2.698568048523214, 0.2584184617970919, 3.950713976768498)
    25 = @(0.05652386341895976, 0.08406217762149737, 0.1139652353035991, 4.233988375685414, 3.845351347943165, 2.499042109264394, 0.2430542141782297, 3.331249627311138)
}

$cols = @("C", "D", "E", "F", "G", "H", "J", "N")

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

Write-Output "Updated $($newValues.Keys.Count) rows"
